$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column I ("Location of data collection") rows 8-23 and 27 held lowercase "de";
# change them to uppercase "DE" (rows 2-7 keep lowercase "de" unchanged).
$deRows = 8..23 + 27
foreach ($r in $deRows) {
    $ws.Cells.Item($r, 9).Value = "DE"
}

# Rows 24-26 held lowercase "ch"; change them to uppercase "CH".
foreach ($r in 24..26) {
    $ws.Cells.Item($r, 9).Value = "CH"
}

# Move the active selection on Sheet1 from I30 to I31, matching the saved view state.
$ws.Range("I31").Select() | Out-Null
